$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Server"
$ws.Move($null, $wb.Worksheets.Item("Network Device"))
$ws2 = $wb.Worksheets.Item("Server")

$src = $wb.Worksheets.Item("Network Device").Range("A4:M4")
$src.Copy()
$dest = $ws2.Range("A4:M4")
$dest.PasteSpecial(-4122)
